# fix lỗi trong report cơ sở. Thêm cột ghi chú trong báo cáo về chi tiêu
#
# 1) Add a new "Đơn phụ phẫu 2" sheet (same shape as "Đơn phụ phẫu 1")
#    placed right before the "Lương" sheet, with header row + 1 data row
#    + a "Tổng" summary row.
# 2) Append a new data row to "Đơn phụ phẫu 1" (before its "Tổng" row) and
#    refresh the "Tổng" row's count/sum.
# 3) Refresh the computed "Lương" sheet cells that roll the new entries up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert "Đơn phụ phẫu 2" right before "Lương"
#
# Cloning "Đơn phụ phẫu 1" (instead of adding a bare sheet) means the
# "Tổng" summary row's filler cells (C..H) come along as the same
# empty-Text cells the template uses, rather than turning into blank
# cells the way a plain `Value = ""` assignment would.
# ---------------------------------------------------------------------
$wsLuong = $wb.Worksheets.Item("Lương")
$wsPP1Template = $wb.Worksheets.Item("Đơn phụ phẫu 1")
$wsPP1Template.Copy($wsLuong)

$wsPP2 = $wb.Worksheets.Item(3)
$wsPP2.Rows.Item(3).Delete()
$wsPP2.Name = "Đơn phụ phẫu 2"

$wsPP2.Cells.Item(1, 8).Value = "Phụ phẫu 2"
$wsPP2.Cells.Item(1, 9).Value = "Công phụ phẫu 2"

$wsPP2.Cells.Item(2, 1).Value = "HD-LUXURY"
$wsPP2.Cells.Item(2, 2).Value = 627
$wsPP2.Cells.Item(2, 3).Value = "'08-04-2024"
$wsPP2.Cells.Item(2, 4).Value = "SÓC TRĂNG"
$wsPP2.Cells.Item(2, 5).Value = "tạ duy hoàng "
$wsPP2.Cells.Item(2, 6).Value = "Cá nhân"
$wsPP2.Cells.Item(2, 7).Value = "Cắt mí"
$wsPP2.Cells.Item(2, 8).Value = "Kha Như Huỳnh "
$wsPP2.Cells.Item(2, 9).Value = 0

# row 3 ("Tổng") kept its label + blank filler cells from the cloned
# template; only its count/sum columns need refreshing
$wsPP2.Cells.Item(3, 2).Value = 1
$wsPP2.Cells.Item(3, 9).Value = 0

# worksheet handles in this host are position-based, so re-resolve
# "Lương" by name now that a new sheet has been spliced in before it
$wsLuong = $wb.Worksheets.Item("Lương")

# ---------------------------------------------------------------------
# 2. "Đơn phụ phẫu 1": add the new order row, refresh "Tổng" row
# ---------------------------------------------------------------------
$wsPP1 = $wb.Worksheets.Item("Đơn phụ phẫu 1")

# push the existing "Tổng" row (row 4) down to row 5, freeing row 4 for
# the new order while keeping the summary row's (empty-string) cell types
$wsPP1.Rows.Item(4).Insert()

$wsPP1.Cells.Item(4, 1).Value = "HD-LUXURY"
$wsPP1.Cells.Item(4, 2).Value = 625
$wsPP1.Cells.Item(4, 3).Value = "'08-04-2024"
$wsPP1.Cells.Item(4, 4).Value = "SÓC TRĂNG"
$wsPP1.Cells.Item(4, 5).Value = "nguyễn thị mỹ chăm"
$wsPP1.Cells.Item(4, 6).Value = "Cá nhân"
$wsPP1.Cells.Item(4, 7).Value = "Cắt mí"
$wsPP1.Cells.Item(4, 8).Value = "Kha Như Huỳnh "
$wsPP1.Cells.Item(4, 9).Value = 50000

$wsPP1.Cells.Item(5, 2).Value = 3
$wsPP1.Cells.Item(5, 9).Value = 200000

# ---------------------------------------------------------------------
# 3. "Lương": refresh the SÓC TRĂNG roll-up figures
# ---------------------------------------------------------------------
$wsLuong.Cells.Item(22, 2).Value = 5.5
$wsLuong.Cells.Item(23, 2).Value = 192500
$wsLuong.Cells.Item(24, 2).Value = 1080357.142857143
$wsLuong.Cells.Item(29, 2).Value = 200000
$wsLuong.Cells.Item(34, 2).Value = 1472857.142857143
$wsLuong.Cells.Item(35, 2).Value = 1472857.142857143

# restore the original active tab (sheet copy/selection above moved it)
$wb.Worksheets.Item(1).Activate()
